$d = $word.ActiveDocument

$quoteOpen  = [char]8220
$quoteClose = [char]8221

# 1) "The information in the file "fort.4" generated by the mesh converter "
#    -> "In principle, the body mass matrix should be input by the users as it contains also the information for the structures above the water surface. "
$find1 = "The information in the file " + $quoteOpen + "fort.4" + $quoteClose + " generated by the mesh converter "
$repl1 = "In principle, the body mass matrix should be input by the users as it contains also the information for the structures above the water surface. "
$d.Content.Find.Execute($find1, $true, $false, $false, $false, $false, $true, 1, $false, $repl1, 2)

# 2) "n.exe can be used for the " -> "n.exe can generate a "
$find2 = "n.exe can be used for the "
$repl2 = "n.exe can generate a "
$d.Content.Find.Execute($find2, $true, $false, $false, $false, $false, $true, 1, $false, $repl2, 2)

# 3) " file. However, the " -> " file in which the body mass matrix is calculated using only the simple information of the wetted body mesh. Therefore, the author should use this body mass matrix with caution. The "
$find3 = " file. However, the "
$repl3 = " file in which the body mass matrix is calculated using only the simple information of the wetted body mesh. Therefore, the author should use this body mass matrix with caution. The "
$d.Content.Find.Execute($find3, $true, $false, $false, $false, $false, $true, 1, $false, $repl3, 2)

# 4) " should be set by the users." -> " should be input by the users."
$find4 = " should be set by the users."
$repl4 = " should be input by the users."
$d.Content.Find.Execute($find4, $true, $false, $false, $false, $false, $true, 1, $false, $repl4, 2)
